$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.627.22"
$ws.Range("E2").Value = "  -1.21%  "

$ws.Range("D3").Value = "3.763.78"
$ws.Range("E3").Value = "  -1.92%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "596.53"
$c.ClearFormats()
$ws.Range("E5").Value = "  -0.86%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "168.07"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "3.763.58"
$ws.Range("E7").Value = "  -1.90%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("E10").Value = "  -1.26%  "

$ws.Range("E11").Value = "  +1.26%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.452"
$c.ClearFormats()
$ws.Range("E12").Value = "  -1.33%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000277"
$c.ClearFormats()
$ws.Range("E13").Value = "  +0.64%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "36.16"
$c.ClearFormats()
$ws.Range("E14").Value = "  -1.79%  "

$ws.Range("D15").Value = "4.396.77"
$ws.Range("E15").Value = "  -1.88%  "

$ws.Range("D16").Value = "3.759.83"
$ws.Range("E16").Value = "  -2.04%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "18.66"
$c.ClearFormats()
$ws.Range("E17").Value = "  +1.90%  "

$ws.Range("D18").Value = "67.627.92"
$ws.Range("E18").Value = "  -1.11%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.19"
$c.ClearFormats()
$ws.Range("E19").Value = "  -1.82%  "

$ws.Range("E20").Value = "  +0.90%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.60"
$c.ClearFormats()
$ws.Range("E21").Value = "  -3.21%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "465.72"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("E23").Value = "  -1.19%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.0000149"
$c.ClearFormats()
$ws.Range("E24").Value = "  -7.14%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "83.56"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("E26").Value = "  -0.84%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "12.06"
$c.ClearFormats()
$ws.Range("E27").Value = "  +0.20%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.46"
$c.ClearFormats()
$ws.Range("E28").Value = "  +1.31%  "

$ws.Range("E29").Value = "  +0.08%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.91"
$c.ClearFormats()
$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("D31").Value = "3.913.33"
$ws.Range("E31").Value = "  -1.87%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.57"
$c.ClearFormats()
$ws.Range("E32").Value = "  -0.98%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "30.36"
$c.ClearFormats()
$ws.Range("E33").Value = "  -2.64%  "

$ws.Range("E34").Value = "  -2.79%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "9.14"
$c.ClearFormats()
$ws.Range("E35").Value = "  -1.45%  "

$ws.Range("D36").Value = "3.724.91"
$ws.Range("E36").Value = "  -2.12%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.ClearFormats()
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.70"
$c.ClearFormats()
$ws.Range("E38").Value = "  -0.27%  "

$ws.Range("E39").Value = "  -0.73%  "

$ws.Range("E40").Value = "  -1.31%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "5.80"
$c.ClearFormats()
$ws.Range("E41").Value = "  -1.42%  "

$ws.Range("E42").Value = "  +0.03%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.311"
$c.ClearFormats()
$ws.Range("E43").Value = "  -0.77%  "

$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("E45").Value = "  -0.18%  "

$ws.Range("E46").Value = "  -1.88%  "

$ws.Range("E47").Value = "  -2.69%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "399.29"
$c.ClearFormats()
$ws.Range("E48").Value = "  -4.78%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.000276"
$c.ClearFormats()
$ws.Range("E49").Value = "  -4.28%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "140.80"
$c.ClearFormats()
$ws.Range("E50").Value = "  -0.42%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "39.47"
$c.ClearFormats()
$ws.Range("E51").Value = "  +6.63%  "
